# Implement database schema migration
# Append one new data row (row 90) to each of the four worksheets,
# mirroring the existing row layout (time, lengths, checksum, decoded values).

$wb = $excel.ActiveWorkbook

$rowsToAdd = @{
    1 = @{
        A = 45769.43330047453
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x5a"
        E = "0xd"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 346
        I = 13
    }
    2 = @{
        A = 45769.2873221875
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x5a"
        E = "0xe"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 346
        I = 14
    }
    3 = @{
        A = 45769.43556526621
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x5a"
        E = "0x3"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 346
        I = 3
    }
    4 = @{
        A = 45769.50058841435
        B = "0x01,0x90"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x01,0x5a"
        E = "0x3"
        F = 400
        G = [double]"9.85046333984776e+23"
        H = 346
        I = 3
    }
}

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $row = $rowsToAdd[$i]
    $newRowIndex = 90

    $ws.Cells.Item($newRowIndex, 1).Value = $row.A
    $ws.Cells.Item($newRowIndex, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($newRowIndex, 2).Value = $row.B
    $ws.Cells.Item($newRowIndex, 3).Value = $row.C
    $ws.Cells.Item($newRowIndex, 4).Value = $row.D
    $ws.Cells.Item($newRowIndex, 5).Value = $row.E

    $ws.Cells.Item($newRowIndex, 6).Value = $row.F
    $ws.Cells.Item($newRowIndex, 7).Value = $row.G
    $ws.Cells.Item($newRowIndex, 8).Value = $row.H
    $ws.Cells.Item($newRowIndex, 9).Value = $row.I
}
